# xldownload for both blockdetails and unitdetails updated
# Rewrites the Block Details sheet: drops the old billing/invoice columns
# (H:O), renames/retypes the remaining headers, replaces the 4 sample
# block rows with 6 new ones, and rewires the manager-email hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the old hyperlinks before we touch the cells they sit on ---
$ws.Hyperlinks.Delete()

# --- new header row -----------------------------------------------------
$ws.Range("A1").Value = "Sno"
$ws.Range("B1").Value = "blockname"
$ws.Range("C1").Value = "blocktype"
$ws.Range("D1").Value = "units"
$ws.Range("E1").Value = "managername"
$ws.Range("F1").Value = "managermobileno"
$ws.Range("G1").Value = "manageremailid"

# --- drop the old FlatRateValue..StartsFrom columns entirely ------------
$ws.Range("H1:O5").Clear()

# --- new block rows -------------------------------------------------------
$blocks = @(
    @(1, "BLOCK-A", "RESIDENTIAL", 5, "SAM",    9999999999, "sam007@gmail.com"),
    @(2, "BLOCK-B", "COMMERCIAL",  6, "JPHN",   8888888888, "john@gmail.com"),
    @(3, "BLOCK-C", "RESIDENTIAL", 4, "SATYA",  7777777777, "sam007@gmail.com"),
    @(4, "BLOCK-D", "COMMERCIAL",  3, "SILVER", 6666666666, "john@gmail.com"),
    @(5, "BLOCK-E", "RESIDENTIAL", 2, "ANU",    5555555555, "sam007@gmail.com"),
    @(6, "BLOCK-F", "COMMERCIAL",  1, "DADY",   4444444444, "john@gmail.com")
)

$row = 2
foreach ($b in $blocks) {
    $ws.Range("A" + $row).Value = $b[0]
    $ws.Range("B" + $row).Value = $b[1]
    $ws.Range("C" + $row).Value = $b[2]
    $ws.Range("D" + $row).Value = $b[3]
    $ws.Range("E" + $row).Value = $b[4]
    $ws.Range("F" + $row).Value = $b[5]
    $ws.Range("G" + $row).Value = $b[6]
    $row = $row + 1
}

# --- reapply the banded row formatting (odd rows mirror row 2, even rows
#     mirror row 3) now that rows 6 & 7 exist ----------------------------
$ws.Range("A2:G2").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$ws.Range("A6:G6").PasteSpecial(-4122)

$ws.Range("A3:G3").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$ws.Range("A7:G7").PasteSpecial(-4122)

# --- manager e-mail hyperlinks ------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:sam007@gmail.com", "", "", "sam007@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:john@gmail.com", "", "", "john@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:sam007@gmail.com", "", "", "sam007@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:john@gmail.com", "", "", "john@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G6"), "mailto:sam007@gmail.com", "", "", "sam007@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G7"), "mailto:john@gmail.com", "", "", "john@gmail.com")
